# LoginHRM.xlsx edit script
# Implements: Tasks sheet cleanup (drop last 2 rows, fix VLOOKUP ranges,
# update Project references for rows 7-11, reuse Projects header labels),
# switch the active tab to "Projects", and tidy up the selections that
# Excel persisted on the Tasks/Projects sheets.

$wb = $excel.ActiveWorkbook

$wsLogin    = $wb.Worksheets.Item("Login")
$wsClients  = $wb.Worksheets.Item("Clients")
$wsProjects = $wb.Worksheets.Item("Projects")
$wsTasks    = $wb.Worksheets.Item("Tasks")

# --- Tasks sheet: header row now reuses the Projects-style labels
# (the dedicated *_TASK shared strings are retired) ---
$wsTasks.Range("A1").Value = "TITLE"
$wsTasks.Range("B1").Value = "START_DATE"
$wsTasks.Range("C1").Value = "END_DATE"
$wsTasks.Range("D1").Value = "HOUR"
$wsTasks.Range("F1").Value = "SUMMARY"
$wsTasks.Range("G1").Value = "DESCRIPTION"
$wsTasks.Range("I1").Value = "CLIENT"

# --- Tasks sheet: widen the VLOOKUP ranges now that Projects has more rows ---
$wsTasks.Range("I2").Formula = "=(VLOOKUP(E2,Projects!A2:B11,2,1))"
$wsTasks.Range("I3").Formula = "=(VLOOKUP(E3,Projects!A3:B7,2,1))"
$wsTasks.Range("I4").Formula = "=(VLOOKUP(E4,Projects!A4:B8,2,1))"
$wsTasks.Range("I5").Formula = "=(VLOOKUP(E5,Projects!A5:B9,2,1))"
$wsTasks.Range("I6").Formula = "=(VLOOKUP(E6,Projects!A6:B10,2,1))"

# --- Tasks sheet rows 7-11: point each task at its own (later) project,
# clear the now-stale explicit style, and fix up the VLOOKUP accordingly ---
$wsTasks.Range("E7").Value = "Project F"
$wsTasks.Range("E7").Style = "Normal"
$wsTasks.Range("I7").Formula = "=(VLOOKUP(E7,Projects!A7:B11,2,1))"

$wsTasks.Range("E8").Value = "Project G"
$wsTasks.Range("E8").Style = "Normal"
$wsTasks.Range("I8").Formula = "=(VLOOKUP(E8,Projects!A8:B11,2,1))"

$wsTasks.Range("E9").Value = "Project H"
$wsTasks.Range("E9").Style = "Normal"
$wsTasks.Range("I9").Formula = "=(VLOOKUP(E9,Projects!A9:B11,2,1))"

$wsTasks.Range("E10").Value = "Project I"
$wsTasks.Range("E10").Style = "Normal"
$wsTasks.Range("I10").Formula = "=(VLOOKUP(E10,Projects!A10:B11,2,1))"

$wsTasks.Range("E11").Value = "Project J"
$wsTasks.Range("E11").Style = "Normal"
$wsTasks.Range("I11").Formula = "=(VLOOKUP(E11,Projects!A11:B11,2,1))"

# --- Tasks sheet: drop the last two (now unused) task rows ---
$wsTasks.Rows.Item(13).Delete()
$wsTasks.Rows.Item(12).Delete()

# --- Tasks sheet: tidy the lingering full-column selection (no longer
# overlapping the now-removed N column) without disturbing the active tab ---
$wsTasks.Range("L1:M1048576").Select()

# --- Projects sheet becomes the active tab/selection ---
$wsProjects.Range("J1:K1048576").Select()
$wsProjects.Activate()

$wb.Save()
